$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "MBA" expense (D7) was paid off / removed from the pending list this
# quarter -- clear its amount and mark it with the same green "paid"
# highlight used elsewhere in the sheet (E5, C6), while keeping its
# existing 2-decimal number format.
$d7 = $ws.Range("D7")
[void]$d7.ClearContents()
$d7.Interior.Color = 5287936
$d7.NumberFormat = "0.00"

# Move the active selection to D10 (single cell).
[void]$ws.Range("D10").Select()
